$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the two extra movies to Sheet1
$ws1.Range("A7").Value = "We are The Millers"
$ws1.Range("A8").Value = "Into The Wild"

# Add the new "Recommended" sheet right after Sheet1
$wsRec = $wb.Worksheets.Add($null, $ws1)
$wsRec.Name = "Recommended"

$wsRec.Range("A1").Value = "Indiana Jones"
$wsRec.Range("A2").Value = "Catch Me If You Can"
$wsRec.Range("A3").Value = "Titanic"
$wsRec.Range("A4").Value = "The Prestige"
$wsRec.Range("A5").Value = "Titanic"
